$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.410.09'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '2.983.05'
$ws.Range('E3').Value = '  +1.54%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '381.32'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.01%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '103.22'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.57%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.546'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('E8').Value = '  +0.00%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.593'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.76%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '36.67'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.47%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.138'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.78%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0861'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '3.451.07'
$ws.Range('E13').Value = '  +1.42%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '7.81'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +3.94%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '18.42'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +2.54%  '
$ws.Range('D16').Value = '2.977.93'
$ws.Range('E16').Value = '  +1.51%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '11.27'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.09%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.67%  '
$ws.Range('D19').Value = '51.398.35'
$ws.Range('E19').Value = '  +0.82%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '3.13'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.46%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '12.62'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('D22').Value = '0.0₃0962'
$ws.Range('E22').Value = '  +0.62%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '70.29'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +2.39%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '267.48'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.14%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '3.23'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +3.50%  '
$ws.Range('E26').Value = '  -3.09%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.50'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.41%  '
$ws.Range('E28').Value = '  -0.07%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '26.08'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('E31').Value = '  -1.39%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '10.33'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.19%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '34.77'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +4.72%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '51.49'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.63%  '
$ws.Range('E35').Value = '  +1.37%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0441'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  +2.88%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '16.84'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +3.24%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.117'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.56'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +3.12%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.84'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.19%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '124.53'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +3.91%  '
$ws.Range('E44').Value = '  +10.37%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '21.61'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.57%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.272'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.37'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +3.23%  '
$ws.Range('D49').Value = '2.033.97'
$ws.Range('E49').Value = '  +2.25%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0331'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.71%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.534'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +16.33%  '
